$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Cells.Item(2, 3).Value = 87
$ws.Cells.Item(2, 5).Value = 0.04316724116202976
$ws.Cells.Item(2, 7).Value = 12
$ws.Cells.Item(3, 3).Value = 88
$ws.Cells.Item(3, 5).Value = 0.007909671282679878
$ws.Cells.Item(3, 7).Value = 8
$ws.Cells.Item(4, 3).Value = 87
$ws.Cells.Item(4, 5).Value = 0.01136965457601704
$ws.Cells.Item(4, 7).Value = 15
$ws.Cells.Item(5, 3).Value = 87
$ws.Cells.Item(5, 5).Value = 0.02256746400640139
$ws.Cells.Item(5, 7).Value = 14
$ws.Cells.Item(6, 3).Value = 87
$ws.Cells.Item(6, 5).Value = 0.01717882544885584
$ws.Cells.Item(6, 7).Value = 15

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Cells.Item(2, 3).Value = 87
$ws.Cells.Item(2, 5).Value = 0.04316724116202976
$ws.Cells.Item(2, 7).Value = 22
$ws.Cells.Item(3, 3).Value = 88
$ws.Cells.Item(3, 5).Value = 0.007909671282679878
$ws.Cells.Item(3, 7).Value = 22
$ws.Cells.Item(4, 3).Value = 87
$ws.Cells.Item(4, 5).Value = 0.01136965457601704
$ws.Cells.Item(4, 7).Value = 29
$ws.Cells.Item(5, 3).Value = 87
$ws.Cells.Item(5, 5).Value = 0.02256746400640139
$ws.Cells.Item(5, 7).Value = 28
$ws.Cells.Item(6, 3).Value = 87
$ws.Cells.Item(6, 5).Value = 0.01717882544885584
$ws.Cells.Item(6, 7).Value = 28

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Cells.Item(2, 3).Value = 87
$ws.Cells.Item(2, 5).Value = 0.04316724116202976
$ws.Cells.Item(2, 7).Value = 43
$ws.Cells.Item(3, 3).Value = 88
$ws.Cells.Item(3, 5).Value = 0.007909671282679878
$ws.Cells.Item(3, 7).Value = 41
$ws.Cells.Item(4, 3).Value = 87
$ws.Cells.Item(4, 5).Value = 0.01136965457601704
$ws.Cells.Item(4, 7).Value = 47
$ws.Cells.Item(5, 3).Value = 87
$ws.Cells.Item(5, 5).Value = 0.02256746400640139
$ws.Cells.Item(5, 7).Value = 51
$ws.Cells.Item(6, 3).Value = 87
$ws.Cells.Item(6, 5).Value = 0.01717882544885584
$ws.Cells.Item(6, 7).Value = 46

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Cells.Item(2, 3).Value = 87
$ws.Cells.Item(2, 5).Value = 0.04316724116202976
$ws.Cells.Item(2, 7).Value = 80
$ws.Cells.Item(3, 3).Value = 88
$ws.Cells.Item(3, 5).Value = 0.007909671282679878
$ws.Cells.Item(3, 7).Value = 64
$ws.Cells.Item(4, 3).Value = 87
$ws.Cells.Item(4, 5).Value = 0.01136965457601704
$ws.Cells.Item(4, 7).Value = 83
$ws.Cells.Item(5, 3).Value = 87
$ws.Cells.Item(5, 5).Value = 0.02256746400640139
$ws.Cells.Item(5, 7).Value = 81
$ws.Cells.Item(6, 3).Value = 87
$ws.Cells.Item(6, 5).Value = 0.01717882544885584
$ws.Cells.Item(6, 7).Value = 68
